# Refined metadata to be additional tab
$wb2 = $excel.ActiveWorkbook
$data = $wb2.Worksheets.Item("data")

# --- Refresh time_taken (F2:F7) on the "data" sheet ---------------------
$data.Range("F2").Value = "2021-10-05 14:19:09.007119"
$data.Range("F3").Value = "2021-10-05 14:19:09.007127"
$data.Range("F4").Value = "2021-10-05 14:19:09.007131"
$data.Range("F5").Value = "2021-10-05 14:19:09.007134"
$data.Range("F6").Value = "2021-10-05 14:19:09.007137"
$data.Range("F7").Value = "2021-10-05 14:19:09.007140"

# --- Add a new "metadata" worksheet, positioned after "data" ------------
$meta = $wb2.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (B1:G1) - same bold/bordered style used by the "data" sheet's header
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats
$data.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)     # xlPasteFormats

# Data row (row 2)
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Aniridia"
$meta.Range("C2").Value = 510

# data_version "2.14" must stay textual (matches the source data, which
# stores it as a string, not a number) -- format the cell as Text first
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.14"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2021-04-14T10:29:35.606157Z"
$meta.Range("F2").Value = "2021-10-05 14:19:09.003747"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/510/?format=json"

# Match the "data" sheet's A2 style (bold/bordered header style) on A2 of metadata
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$meta.Range("A1").Select() | Out-Null
